$wb = $excel.ActiveWorkbook

# --- Overview sheet: update status text (shared by E/F columns, rows 2-3) ---
$ovw = $wb.Worksheets.Item("Overview")
$ovw.Range("E2:F3").Value = "Handed back: in sync with en-US"

# --- zh-cn sheet ---
$zh = $wb.Worksheets.Item("zh-cn")
$zh.Range("C2:C3").Value = "Handed back: in sync with en-US"

$zh.Range("I2").Value = "ae028673-f9a0-4771-8ab2-37d5a1c491fe.md"
$zh.Range("J2").Value = "ae028673-f9a0-4771-8ab2-37d5a1c491fe.e8be8300cd87d911dc40d54f726a2dda27c158dc.zh-cn.xlf"
$zh.Range("K2").Value = "2016-08-28 06:52:32"

$zh.Range("I3").Value = "cd54ed06-4bee-4486-a1f0-1dc02011ca95.md"
$zh.Range("J3").Value = "cd54ed06-4bee-4486-a1f0-1dc02011ca95.12ffa7c52420325959c5e575bf27d9b8c17ed3d5.zh-cn.xlf"
$zh.Range("K3").Value = "2016-08-28 06:52:32"

# --- de-de sheet ---
$de = $wb.Worksheets.Item("de-de")
$de.Range("C2:C3").Value = "Handed back: in sync with en-US"

$de.Range("I2").Value = "ae028673-f9a0-4771-8ab2-37d5a1c491fe.md"
$de.Range("J2").Value = "ae028673-f9a0-4771-8ab2-37d5a1c491fe.e8be8300cd87d911dc40d54f726a2dda27c158dc.de-de.xlf"
$de.Range("K2").Value = "2016-08-28 06:52:38"

$de.Range("I3").Value = "cd54ed06-4bee-4486-a1f0-1dc02011ca95.md"
$de.Range("J3").Value = "cd54ed06-4bee-4486-a1f0-1dc02011ca95.12ffa7c52420325959c5e575bf27d9b8c17ed3d5.de-de.xlf"
$de.Range("K3").Value = "2016-08-28 06:52:38"
